# Update the workbook for the 2022-12-12 data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the header label to reflect the new "through" date.
$ws.Name = "Through 2022-12-12"
$ws.Range("I1").Value = "2022 (through 12-12)"

# Update the December row's 2022 value, and the Total row's 2022 value.
$ws.Range("I13").Value = 53
$ws.Range("I14").Value = 1569
